$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "4111-03550-C"
$ws.Range("B2").Value = "Besar"
$ws.Range("C2").Value = 14

$ws.Range("A3").Value = "4111-03550-C"
$ws.Range("B3").Value = "Besar"
$ws.Range("C3").Value = 14
